# Weekly update: insert a new week's worth of Chirimoya price rows
# (Especial / Primera / Segunda, Cultivar IV Región, date 2021-10-13)
# right before the existing row 44, shifting all following rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at position 44 (pushes old rows 44.. down to 47..)
$ws.Rows("44:46").Insert()

# Shared values across the three new rows
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$fecha       = 44482
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100107
$producto    = "Otros"
$categoriaId = 100107002
$categoria   = "Chirimoya"
$variedad    = "Cultivar IV Región"
$unidad      = "`$/kilo (en caja de 15 kilos)"
$origen      = "Provincia de Limarí"
$kgUnidad    = 1

$newRows = @(
    @{ Row = 44; Calidad = "Especial"; Volumen = 170; PMin = 2700; PMax = 2700; PProm = 2700; PKg = 2700 },
    @{ Row = 45; Calidad = "Primera";  Volumen = 225; PMin = 2400; PMax = 2400; PProm = 2400; PKg = 2400 },
    @{ Row = 46; Calidad = "Segunda";  Volumen = 200; PMin = 2000; PMax = 2000; PProm = 2000; PKg = 2000 }
)

foreach ($d in $newRows) {
    $r = $d.Row
    $ws.Range("A$r").Value = $mercadoId
    $ws.Range("B$r").Value = $mercado
    $ws.Range("C$r").Value = $region
    $ws.Range("D$r").Value = $fecha
    $ws.Range("E$r").Value = $codreg
    $ws.Range("F$r").Value = $tipo
    $ws.Range("G$r").Value = $productoId
    $ws.Range("H$r").Value = $producto
    $ws.Range("I$r").Value = $categoriaId
    $ws.Range("J$r").Value = $categoria
    $ws.Range("K$r").Value = $variedad
    $ws.Range("L$r").Value = $d.Calidad
    $ws.Range("M$r").Value = $d.Volumen
    $ws.Range("N$r").Value = $d.PMin
    $ws.Range("O$r").Value = $d.PMax
    $ws.Range("P$r").Value = $d.PProm
    $ws.Range("Q$r").Value = $unidad
    $ws.Range("R$r").Value = $origen
    $ws.Range("S$r").Value = $d.PKg
    $ws.Range("T$r").Value = $kgUnidad
}
